$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.0434
$ws.Range("A12").Value = -21.59440000000001
$ws.Range("E14").Value = 16.43170000000001
$ws.Range("E26").Value = 16.367
$ws.Range("A27").Value = -21.78429999999998
$ws.Range("E31").Value = 16.4695
$ws.Range("A32").Value = -21.2518
$ws.Range("E35").Value = 16.4158
$ws.Range("A36").Value = -20.3038
$ws.Range("E37").Value = 16.65820000000001
$ws.Range("A38").Value = -19.4251
$ws.Range("E45").Value = 16.5457
$ws.Range("A46").Value = -21.81180000000001
$ws.Range("E52").Value = 17.1962
$ws.Range("A54").Value = -21.60339999999999
$ws.Range("A55").Value = -22.5005
$ws.Range("A56").Value = -22.2113
$ws.Range("E57").Value = 16.65060000000001
$ws.Range("A67").Value = -21.43699999999997
$ws.Range("A69").Value = -21.57319999999996
$ws.Range("A72").Value = -21.82079999999999
$ws.Range("E81").Value = 16.41189999999999
$ws.Range("A83").Value = -21.9218
$ws.Range("E83").Value = 16.4559
$ws.Range("A86").Value = -22.02470000000001
$ws.Range("A91").Value = -21.40880000000001
$ws.Range("A93").Value = -21.334
$ws.Range("A99").Value = -19.97099999999999
$ws.Range("E100").Value = 16.302
$ws.Range("E102").Value = 16.71609999999999
